# Generate Report for Handback
# Update the "generated at" timestamps recorded on the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on the Overview sheet, and the matching
# "Correspond Handoff Datetime" on the de-de detail sheet share the same
# underlying timestamp value - update both occurrences together.
$wsOverview.Range("G2").Value = "2016-09-02 09:15:58"
$wsDeDe.Range("H2").Value = "2016-09-02 09:15:58"

# zh-cn detail sheet: handoff + handback timestamps for the first row.
$wsZhCn.Range("H2").Value = "2016-09-02 09:15:53"
$wsZhCn.Range("K2").Value = "2016-09-02 09:16:30"

# de-de detail sheet: handback timestamp for the first row.
$wsDeDe.Range("K2").Value = "2016-09-02 09:16:38"
